$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 16 de Septiembre de 2020 a las 21:00'
$ws.Cells.Item(4, 2).Value = 6803927
$ws.Cells.Item(4, 3).Value = 15780
$ws.Cells.Item(4, 4).Value = 4088443
$ws.Cells.Item(4, 5).Value = 2514631
$ws.Cells.Item(4, 7).Value = 656
$ws.Cells.Item(4, 8).Value = 200853
$ws.Cells.Item(16, 2).Value = 404888
$ws.Cells.Item(16, 3).Value = 9784
$ws.Cells.Item(16, 4).Value = 90335
$ws.Cells.Item(16, 5).Value = 283508
$ws.Cells.Item(16, 7).Value = 46
$ws.Cells.Item(16, 8).Value = 31045
$ws.Cells.Item(51, 1).Value = 'Etiopia'
$ws.Cells.Item(51, 2).Value = 66224
$ws.Cells.Item(51, 3).Value = 738
$ws.Cells.Item(51, 4).Value = 26665
$ws.Cells.Item(51, 5).Value = 38514
$ws.Cells.Item(51, 7).Value = 10
$ws.Cells.Item(51, 8).Value = 1045
$ws.Cells.Item(52, 1).Value = 'Portugal'
$ws.Cells.Item(52, 2).Value = 65626
$ws.Cells.Item(52, 3).Value = 605
$ws.Cells.Item(52, 4).Value = 44528
$ws.Cells.Item(52, 5).Value = 19220
$ws.Cells.Item(52, 7).Value = 3
$ws.Cells.Item(52, 8).Value = 1878
$ws.Cells.Item(73, 2).Value = 31799
$ws.Cells.Item(73, 3).Value = 250
$ws.Cells.Item(73, 5).Value = 6647
$ws.Cells.Item(73, 7).Value = 1
$ws.Cells.Item(73, 8).Value = 1788
$ws.Cells.Item(99, 2).Value = 9595
$ws.Cells.Item(99, 3).Value = 17
$ws.Cells.Item(99, 4).Value = 9235
$ws.Cells.Item(99, 5).Value = 295
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 65
$ws.Cells.Item(121, 2).Value = 4876
$ws.Cells.Item(121, 3).Value = 73
$ws.Cells.Item(121, 4).Value = 4131
$ws.Cells.Item(121, 5).Value = 636
$ws.Cells.Item(121, 7).Value = 1
$ws.Cells.Item(121, 8).Value = 109
$ws.Cells.Item(128, 1).Value = 'Siria'
$ws.Cells.Item(128, 2).Value = 3654
$ws.Cells.Item(128, 3).Value = 40
$ws.Cells.Item(128, 4).Value = 889
$ws.Cells.Item(128, 5).Value = 2602
$ws.Cells.Item(128, 7).Value = 3
$ws.Cells.Item(128, 8).Value = 163
$ws.Cells.Item(129, 1).Value = 'Birmania'
$ws.Cells.Item(129, 2).Value = 3636
$ws.Cells.Item(129, 3).Value = 134
$ws.Cells.Item(129, 4).Value = 832
$ws.Cells.Item(129, 5).Value = 2765
$ws.Cells.Item(129, 7).Value = 4
$ws.Cells.Item(129, 8).Value = 39
$ws.Cells.Item(133, 2).Value = 3440
$ws.Cells.Item(133, 3).Value = 12
$ws.Cells.Item(133, 4).Value = 1851
$ws.Cells.Item(133, 5).Value = 1482
$ws.Cells.Item(133, 7).Value = 2
$ws.Cells.Item(133, 8).Value = 107
$ws.Cells.Item(134, 1).Value = 'Guadalupe'
$ws.Cells.Item(134, 2).Value = 3426
$ws.Cells.Item(134, 3).Value = 346
$ws.Cells.Item(134, 4).Value = 837
$ws.Cells.Item(134, 5).Value = 2563
$ws.Cells.Item(134, 7).Value = 2
$ws.Cells.Item(134, 8).Value = 26
$ws.Cells.Item(135, 1).Value = 'Somalia'
$ws.Cells.Item(135, 2).Value = 3390
$ws.Cells.Item(135, 4).Value = 2812
$ws.Cells.Item(135, 5).Value = 480
$ws.Cells.Item(135, 8).Value = 98
$ws.Cells.Item(136, 1).Value = 'Mayotte'
$ws.Cells.Item(136, 2).Value = 3374
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(136, 4).Value = 2964
$ws.Cells.Item(136, 5).Value = 370
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 40
$ws.Cells.Item(137, 1).Value = 'Trinidad yTobago'
$ws.Cells.Item(137, 2).Value = 3293
$ws.Cells.Item(137, 3).Value = 70
$ws.Cells.Item(137, 4).Value = 810
$ws.Cells.Item(137, 5).Value = 2426
$ws.Cells.Item(137, 7).Value = 1
$ws.Cells.Item(137, 8).Value = 57
$ws.Cells.Item(138, 1).Value = 'Sri Lanka'
$ws.Cells.Item(138, 2).Value = 3271
$ws.Cells.Item(138, 4).Value = 3021
$ws.Cells.Item(138, 5).Value = 237
$ws.Cells.Item(138, 8).Value = 13
$ws.Cells.Item(139, 1).Value = 'Aruba'
$ws.Cells.Item(139, 2).Value = 3152
$ws.Cells.Item(139, 4).Value = 1610
$ws.Cells.Item(139, 5).Value = 1520
$ws.Cells.Item(139, 8).Value = 22
$ws.Cells.Item(191, 2).Value = 181
$ws.Cells.Item(191, 3).Value = 3
$ws.Cells.Item(191, 4).Value = 137
$ws.Cells.Item(191, 5).Value = 43
$ws.Cells.Item(204, 1).Value = 'Santa Lucia'
$ws.Cells.Item(205, 1).Value = 'Timor Oriental'
$ws.Cells.Item(214, 1).Value = 'Montserrat'
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1
$ws.Cells.Item(215, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0
